$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 89, shifting existing rows 89-154 down to 90-155
$ws.Rows("89:89").Insert()

# Populate the newly inserted row 89 with the new record
$ws.Range("A89").Value = 11
$ws.Range("B89").Value = "Vega Monumental Concepción"
$ws.Range("C89").Value = "Bíobío"
$ws.Range("D89").Value = 45001
$ws.Range("E89").Value = 8
$ws.Range("F89").Value = "Fruta"
$ws.Range("G89").Value = 100108
$ws.Range("H89").Value = "Tropicales y subtropicales"
$ws.Range("I89").Value = 100108002
$ws.Range("J89").Value = "Mango"
$ws.Range("K89").Value = "Sin especificar"
$ws.Range("L89").Value = "Primera"
$ws.Range("M89").Value = 200
$ws.Range("N89").Value = 7000
$ws.Range("O89").Value = 7500
$ws.Range("P89").Value = 7250
$ws.Range("Q89").Value = "$/bandeja 4 kilos"
$ws.Range("R89").Value = "Perú"
$ws.Range("S89").Value = 1812
$ws.Range("T89").Value = 4
